# Update the data dictionary: remove the "Creditos_Periodo" and
# "Creditos_Acumulados" fields (columns N:O) from the worksheet.
# Deleting these two whole columns shifts every later column left by two,
# which is what the published dictionary update does (headers/values for
# Contracreditos_*, Apropiacion_Definitiva, Aplazamientos, ... all move
# left to close the gap, and the now-unused shared strings drop out of the
# table automatically on save).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1:O2").EntireColumn.Delete() | Out-Null

# Restore the selection to match the saved state of the workbook.
$ws.Range("Z1").Select() | Out-Null
